$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-29 Friday" "2024-03-30 Saturday"

Replace-Text "200÷2=100, 0" "891÷6=148, 3"
Replace-Text "567÷8=70, 7" "131÷2=65, 1"
Replace-Text "496÷5=99, 1" "728÷3=242, 2"
Replace-Text "948÷3=316, 0" "125÷7=17, 6"
Replace-Text "723÷5=144, 3" "392÷8=49, 0"

Replace-Text "749÷3=249, 2" "659÷7=94, 1"
Replace-Text "560÷2=280, 0" "416÷8=52, 0"
Replace-Text "957÷3=319, 0" "270÷4=67, 2"
Replace-Text "982÷4=245, 2" "891÷4=222, 3"
Replace-Text "266÷6=44, 2" "443÷4=110, 3"

Replace-Text "662÷4=165, 2" "164÷3=54, 2"
Replace-Text "375÷5=75, 0" "106÷8=13, 2"
Replace-Text "545÷2=272, 1" "500÷9=55, 5"
Replace-Text "984÷8=123, 0" "797÷4=199, 1"
Replace-Text "350÷8=43, 6" "634÷6=105, 4"

Replace-Text "739÷4=184, 3" "185÷3=61, 2"
Replace-Text "374÷8=46, 6" "919÷5=183, 4"
Replace-Text "930÷2=465, 0" "985÷5=197, 0"
Replace-Text "687÷3=229, 0" "705÷5=141, 0"
Replace-Text "847÷2=423, 1" "907÷9=100, 7"

Replace-Text "651÷3=217, 0" "588÷8=73, 4"
Replace-Text "218÷9=24, 2" "636÷5=127, 1"
Replace-Text "949÷4=237, 1" "389÷9=43, 2"
Replace-Text "924÷3=308, 0" "274÷9=30, 4"
Replace-Text "985÷4=246, 1" "931÷6=155, 1"

Write-Output "Done"
